$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.394.14'
$ws.Range('E2').Value = '  -1.43%  '

$ws.Range('D3').Value = '1.630.91'
$ws.Range('E3').Value = '  -1.62%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9999'
$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('E5').Value = '  +0.08%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '298.16'
$ws.Range('E6').Value = '  -1.84%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3751'
$ws.Range('E7').Value = '  -1.53%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '50.40'
$ws.Range('E8').Value = '  -1.62%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.3493'
$ws.Range('E9').Value = '  -3.46%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.08034'
$ws.Range('E10').Value = '  -2.18%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '1.205'
$ws.Range('E11').Value = '  -2.86%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.9998'
$ws.Range('E12').Value = '  -0.01%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '21.88'
$ws.Range('E13').Value = '  -3.39%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.309'
$ws.Range('E14').Value = '  -3.30%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.225'
$ws.Range('E15').Value = '  -2.77%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.00001196'
$ws.Range('E16').Value = '  -3.05%  '

$ws.Range('D17').Value = '1.629.40'
$ws.Range('E17').Value = '  -0.86%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '94.82'
$ws.Range('E18').Value = '  -2.84%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06938'
$ws.Range('E19').Value = '  -0.75%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.630'
$ws.Range('E20').Value = '  -2.99%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '17.34'
$ws.Range('E21').Value = '  -1.93%  '

$ws.Range('E22').Value = '  +0.09%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '12.35'
$ws.Range('E23').Value = '  -3.99%  '

$ws.Range('D24').Value = '23.393.89'
$ws.Range('E24').Value = '  -1.44%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.427'
$ws.Range('E25').Value = '  -3.36%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.951'
$ws.Range('E26').Value = '  -3.10%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '20.74'
$ws.Range('E27').Value = '  -2.55%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '151.66'
$ws.Range('E28').Value = '  -0.22%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.182'
$ws.Range('E29').Value = '  -0.49%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '131.76'
$ws.Range('E30').Value = '  -2.07%  '

$ws.Range('D31').Value = '1.807.22'
$ws.Range('E31').Value = '  -1.60%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.782'
$ws.Range('E32').Value = '  -2.74%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '2.113'
$ws.Range('E33').Value = '  -3.24%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '11.06'
$ws.Range('E34').Value = '  -8.15%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.9697'
$ws.Range('E35').Value = '  -8.82%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.02679'

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.08709'
$ws.Range('E37').Value = '  -0.93%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.2424'
$ws.Range('E38').Value = '  -4.06%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.848'
$ws.Range('E39').Value = '  -4.36%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.06713'
$ws.Range('E40').Value = '  -4.87%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '12.70'
$ws.Range('E41').Value = '  -2.98%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.6804'
$ws.Range('E42').Value = '  -3.34%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.286'
$ws.Range('E43').Value = '  -3.74%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '15.25'
$ws.Range('E44').Value = '  -4.99%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.9998'
$ws.Range('E45').Value = '  +0.07%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.6320'
$ws.Range('E46').Value = '  -3.29%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.231'
$ws.Range('E47').Value = '  -4.02%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '3.889'
$ws.Range('E48').Value = '  -2.26%  '

$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '126.67'
$ws.Range('E49').Value = '  -1.12%  '

$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.07643'
$ws.Range('E50').Value = '  -3.84%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.216'
$ws.Range('E51').Value = '  +1.93%  '
